$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.48
$wsSummary.Range("B4").Value = 0.58
$wsSummary.Range("B5").Value = 0.15
$wsSummary.Range("B6").Value = 76
$wsSummary.Range("B7").Value = 38
$wsSummary.Range("B9").Value = 50

# ---------------------------------------------------------------------------
# Strategy Status sheet - HighProbConvergence row (row 3)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C3").Value = 100.12
$wsStatus.Range("D3").Value = 4
$wsStatus.Range("E3").Value = 0.12
$wsStatus.Range("F3").Value = 0.12
$wsStatus.Range("G3").Value = 75

# ---------------------------------------------------------------------------
# All Trades sheet - close trade #77 (row 78)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Range("G78").Value = 0.68
$wsAll.Range("H78").Value = "CLOSED"
$wsAll.Range("I78").Value = 9.6774
$wsAll.Range("J78").Value = 0.06
$wsAll.Range("K78").Value = 100.12
$wsAll.Range("L78").Value = "early_exit"
$wsAll.Range("M78").Value = 0.1

# New row 107 - trade #106 (momentum, DOWN, OPEN)
$wsAll.Range("A107").Value = 106
$wsAll.Range("B107").Value = "'2026-02-18"
$wsAll.Range("C107").Value = "00:20:10"
$wsAll.Range("D107").Value = "momentum"
$wsAll.Range("E107").Value = "DOWN"
$wsAll.Range("F107").Value = 0.62
$wsAll.Range("G107").Value = "'"
$wsAll.Range("H107").Value = "OPEN"
$wsAll.Range("I107").Value = 0
$wsAll.Range("J107").Value = 0
$wsAll.Range("K107").Value = 99.6787371310913
$wsAll.Range("L107").Value = "'"
$wsAll.Range("M107").Value = 0
$wsAll.Range("N107").Value = 0
$wsAll.Range("O107").Value = 0
$wsAll.Range("P107").Value = 0.9
$wsAll.Range("Q107").Value = "Downward momentum: -3.810% over 10 samples"

# ---------------------------------------------------------------------------
# momentum sheet - new row 25 - trade #106
# ---------------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")
$wsMomentum.Range("A25").Value = 106
$wsMomentum.Range("B25").Value = "'2026-02-18"
$wsMomentum.Range("C25").Value = "00:20:10"
$wsMomentum.Range("D25").Value = "momentum"
$wsMomentum.Range("E25").Value = "DOWN"
$wsMomentum.Range("F25").Value = 0.62
$wsMomentum.Range("G25").Value = "'"
$wsMomentum.Range("H25").Value = "OPEN"
$wsMomentum.Range("I25").Value = 0
$wsMomentum.Range("J25").Value = 0
$wsMomentum.Range("K25").Value = 99.6787371310913
$wsMomentum.Range("L25").Value = 0
$wsMomentum.Range("M25").Value = 0
$wsMomentum.Range("N25").Value = 0.9
$wsMomentum.Range("O25").Value = "Downward momentum: -3.810% over 10 samples"
$wsMomentum.Range("P25").Value = "'"
$wsMomentum.Range("Q25").Value = 0

# ---------------------------------------------------------------------------
# HighProbConvergence sheet - close trade #77 (row 5)
# ---------------------------------------------------------------------------
$wsHPC = $wb.Worksheets.Item("HighProbConvergence")
$wsHPC.Range("G5").Value = 0.68
$wsHPC.Range("H5").Value = "CLOSED"
$wsHPC.Range("I5").Value = 9.6774
$wsHPC.Range("J5").Value = 0.06
$wsHPC.Range("K5").Value = 100.12
$wsHPC.Range("P5").Value = "early_exit"
$wsHPC.Range("Q5").Value = 0.1
